# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet a new (blank) column is inserted
# immediately before the existing "Late" column (column N), pushing the
# "Late" / blank-heading / "Outstanding" block one column to the right
# (N->O, O->P, P->Q) and shifting all of their values along with it.
# The sheet also becomes the active/selected sheet/tab, with cell K19
# selected, and the "Input" sheet is no longer the selected tab.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14) -- this shifts the
# "Late", (blank heading) and "Outstanding" columns (and all of their
# row values) one column to the right, exactly like using Excel's
# Insert > Sheet Columns on column N.
$wsSchedule.Columns.Item(14).Insert()

# Match the width of the newly inserted column to its left neighbour
# (column M / "In Advance"), which is what Excel does visually here.
$wsSchedule.Columns.Item(14).ColumnWidth = $wsSchedule.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab (was "Input") and park
# the selection at K19, matching the saved view state in the workbook.
[void]$wsSchedule.Activate()
[void]$wsSchedule.Range("K19").Select()
